# Lotomania_edt.xlsx - append the latest lottery draws (contests 2826-2842)
# to the LOTOMANIA sheet, right below the previous last row (319), then
# leave the freshly-typed block selected like Excel does right after a
# paste of new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one drawn contest: [contest number, 20 sorted balls]
$data = @(
    ,@(2826,9,10,15,16,17,24,25,26,31,32,41,42,45,48,71,72,79,84,87,95)
    ,@(2827,10,11,17,18,19,24,27,31,37,42,43,57,58,64,66,71,85,93,94,95)
    ,@(2828,3,5,10,17,22,23,29,39,48,56,57,63,67,71,74,82,84,87,94,99)
    ,@(2829,2,7,9,16,22,24,27,40,45,57,62,64,65,72,75,86,90,94,96,98)
    ,@(2830,0,12,14,18,22,23,25,29,30,33,34,49,55,67,71,80,83,84,92,93)
    ,@(2831,4,5,7,15,16,24,25,28,38,41,43,47,61,64,69,70,77,81,85,87)
    ,@(2832,3,14,18,20,22,23,25,30,40,47,49,53,63,66,82,87,90,92,93,98)
    ,@(2833,15,19,22,30,33,40,44,46,49,54,57,61,66,70,71,75,81,82,97,99)
    ,@(2834,9,13,27,40,41,42,43,47,53,61,62,65,66,67,72,77,88,92,93,94)
    ,@(2835,7,8,24,27,29,33,40,46,49,50,66,68,75,81,82,86,88,92,98,99)
    ,@(2836,15,20,24,26,27,46,47,50,51,54,55,56,58,59,66,82,85,88,90,99)
    ,@(2837,1,6,8,10,11,15,19,22,38,43,50,61,62,65,69,74,77,84,85,99)
    ,@(2838,5,6,12,16,23,24,31,34,35,36,43,53,58,63,82,85,86,87,97,99)
    ,@(2839,2,9,13,23,26,38,45,48,51,56,61,69,73,76,77,83,87,93,95,98)
    ,@(2840,3,4,7,19,22,24,30,32,38,39,42,43,46,47,50,51,62,68,72,95)
    ,@(2841,0,7,8,10,13,16,24,26,31,37,39,41,42,43,52,57,70,73,75,82)
    ,@(2842,5,6,10,11,12,21,25,27,35,44,48,54,63,64,67,75,76,83,84,86)
)

$startRow = 320
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($col = 0; $col -lt $rowVals.Count; $col++) {
        $ws.Cells.Item($r, $col + 1).Value = $rowVals[$col]
    }
}

# The last 11 rows (326-336) were pasted in a separate step from a
# differently-formatted source, so they pick up their own (still default
# looking) style records instead of reusing the very first one. Reproduce
# that by touching the column-A cells and the B:U cells of that block with
# two distinct, purely cosmetic format tweaks.
$lastBlock = $ws.Range("A326:A336")
$lastBlock.Locked = $true

$lastBlockRest = $ws.Range("B326:U336")
$lastBlockRest.NumberFormat = "General"

# Match the post-paste selection: the freshly added block (minus the
# contest-number column) ends up selected, anchored at its first cell.
[void]$ws.Range("B326:U336").Select()
